# Fruta / hortaliza, semanal
# Insert two new weekly observation pairs into the Chirimoya price history
# sheet: one pair ("Especial"/"Primera", $/bandeja 10 kilos) dated 44846,
# inserted before the existing row with date 44160, and another pair
# ("Especial"/"Primera", $/bandeja 10 kilos) dated 44845, inserted before
# the existing row with date 44454.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 blank rows at row 24 (pushes old rows 24-31 down to 26-33) ---
$ws.Rows("24:25").Insert()

# --- Insert 2 blank rows at row 28 (pushes old rows 26-31, now at 26-31 minus
#     the two already-shifted rows 24/25 kept at 26/27, i.e. rows 28-33, down
#     to 30-35) ---
$ws.Rows("28:29").Insert()

# Common (unchanged across every row of this block) column values
$marketId = 11
$market   = "Vega Monumental Concepción"
$region   = "Bíobío"
$codreg   = 8
$tipo     = "Fruta"
$prodId   = 100107
$prod     = "Otros"
$catId    = 100107002
$cat      = "Chirimoya"
$variedad = "Cultivar IV Región"
$origen   = "Provincia de Limarí"

# --- New row 24 ---
$ws.Cells.Item(24, 1).Value2 = $marketId
$ws.Cells.Item(24, 2).Value2 = $market
$ws.Cells.Item(24, 3).Value2 = $region
$ws.Cells.Item(24, 4).Value2 = 44846
$ws.Cells.Item(24, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
$ws.Cells.Item(24, 5).Value2 = $codreg
$ws.Cells.Item(24, 6).Value2 = $tipo
$ws.Cells.Item(24, 7).Value2 = $prodId
$ws.Cells.Item(24, 8).Value2 = $prod
$ws.Cells.Item(24, 9).Value2 = $catId
$ws.Cells.Item(24, 10).Value2 = $cat
$ws.Cells.Item(24, 11).Value2 = $variedad
$ws.Cells.Item(24, 12).Value2 = "Especial"
$ws.Cells.Item(24, 13).Value2 = 50
$ws.Cells.Item(24, 14).Value2 = 23000
$ws.Cells.Item(24, 15).Value2 = 23000
$ws.Cells.Item(24, 16).Value2 = 23000
$ws.Cells.Item(24, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(24, 18).Value2 = $origen
$ws.Cells.Item(24, 19).Value2 = 2300
$ws.Cells.Item(24, 20).Value2 = 10

# --- New row 25 ---
$ws.Cells.Item(25, 1).Value2 = $marketId
$ws.Cells.Item(25, 2).Value2 = $market
$ws.Cells.Item(25, 3).Value2 = $region
$ws.Cells.Item(25, 4).Value2 = 44846
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
$ws.Cells.Item(25, 5).Value2 = $codreg
$ws.Cells.Item(25, 6).Value2 = $tipo
$ws.Cells.Item(25, 7).Value2 = $prodId
$ws.Cells.Item(25, 8).Value2 = $prod
$ws.Cells.Item(25, 9).Value2 = $catId
$ws.Cells.Item(25, 10).Value2 = $cat
$ws.Cells.Item(25, 11).Value2 = $variedad
$ws.Cells.Item(25, 12).Value2 = "Primera"
$ws.Cells.Item(25, 13).Value2 = 100
$ws.Cells.Item(25, 14).Value2 = 20000
$ws.Cells.Item(25, 15).Value2 = 21000
$ws.Cells.Item(25, 16).Value2 = 20500
$ws.Cells.Item(25, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(25, 18).Value2 = $origen
$ws.Cells.Item(25, 19).Value2 = 2050
$ws.Cells.Item(25, 20).Value2 = 10

# --- New row 28 ---
$ws.Cells.Item(28, 1).Value2 = $marketId
$ws.Cells.Item(28, 2).Value2 = $market
$ws.Cells.Item(28, 3).Value2 = $region
$ws.Cells.Item(28, 4).Value2 = 44845
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
$ws.Cells.Item(28, 5).Value2 = $codreg
$ws.Cells.Item(28, 6).Value2 = $tipo
$ws.Cells.Item(28, 7).Value2 = $prodId
$ws.Cells.Item(28, 8).Value2 = $prod
$ws.Cells.Item(28, 9).Value2 = $catId
$ws.Cells.Item(28, 10).Value2 = $cat
$ws.Cells.Item(28, 11).Value2 = $variedad
$ws.Cells.Item(28, 12).Value2 = "Especial"
$ws.Cells.Item(28, 13).Value2 = 100
$ws.Cells.Item(28, 14).Value2 = 23000
$ws.Cells.Item(28, 15).Value2 = 23000
$ws.Cells.Item(28, 16).Value2 = 23000
$ws.Cells.Item(28, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(28, 18).Value2 = $origen
$ws.Cells.Item(28, 19).Value2 = 2300
$ws.Cells.Item(28, 20).Value2 = 10

# --- New row 29 ---
$ws.Cells.Item(29, 1).Value2 = $marketId
$ws.Cells.Item(29, 2).Value2 = $market
$ws.Cells.Item(29, 3).Value2 = $region
$ws.Cells.Item(29, 4).Value2 = 44845
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
$ws.Cells.Item(29, 5).Value2 = $codreg
$ws.Cells.Item(29, 6).Value2 = $tipo
$ws.Cells.Item(29, 7).Value2 = $prodId
$ws.Cells.Item(29, 8).Value2 = $prod
$ws.Cells.Item(29, 9).Value2 = $catId
$ws.Cells.Item(29, 10).Value2 = $cat
$ws.Cells.Item(29, 11).Value2 = $variedad
$ws.Cells.Item(29, 12).Value2 = "Primera"
$ws.Cells.Item(29, 13).Value2 = 100
$ws.Cells.Item(29, 14).Value2 = 21000
$ws.Cells.Item(29, 15).Value2 = 21000
$ws.Cells.Item(29, 16).Value2 = 21000
$ws.Cells.Item(29, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(29, 18).Value2 = $origen
$ws.Cells.Item(29, 19).Value2 = 2100
$ws.Cells.Item(29, 20).Value2 = 10

Write-Host "Done"
